# feat: add 2022-Q4 data
#
# 1) "总计" (summary) sheet: insert a new row at position 2 for "2022-Q4"
#    and shift the existing rows down.
# 2) Insert a brand-new "2022-Q4" worksheet (copied from the "2022-Q3"
#    sheet so it keeps identical layout/formatting) right after "总计",
#    populated with the 2022-Q4 fund-holding detail rows.

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text even when the string looks
    # numeric (e.g. "005262" or "1.55") so leading zeros / exact string
    # formatting survive, then drop the temporary text number-format so
    # the cell is left with no special style applied.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# The insert shifted the old rows down but column A (the plain 0-based
# row index) needs to stay a simple sequence matching the row position,
# so re-stamp A2:A8 explicitly after the shift.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.12

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2. Create the "2022-Q4" detail sheet (copy formatting from "2022-Q3")
# ---------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Trim the copied sheet down to the 5 data rows we need (rows 7-26 held
# the old 2022-Q3 detail rows we don't want here).
$q4.Rows("7:26").Delete()

# Row 2
$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "005262"
Set-TextValue $q4.Range("C2") "鑫元欣享灵活配置混合A"
Set-TextValue $q4.Range("D2") "1.55"
Set-TextValue $q4.Range("E2") "83.07"
Set-TextValue $q4.Range("F2") "3.42"
Set-TextValue $q4.Range("G2") "0.0530"
$q4.Range("H2").Value = 8

# Row 3
$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "009395"
Set-TextValue $q4.Range("C3") "鑫元安鑫回报混合A"
Set-TextValue $q4.Range("D3") "3.61"
Set-TextValue $q4.Range("E3") "20.13"
Set-TextValue $q4.Range("F3") "0.97"
Set-TextValue $q4.Range("G3") "0.0350"
$q4.Range("H3").Value = 8

# Row 4
$q4.Range("A4").Value = 2
Set-TextValue $q4.Range("B4") "005263"
Set-TextValue $q4.Range("C4") "鑫元欣享灵活配置混合C"
Set-TextValue $q4.Range("D4") "1.00"
Set-TextValue $q4.Range("E4") "83.07"
Set-TextValue $q4.Range("F4") "3.42"
Set-TextValue $q4.Range("G4") "0.0342"
$q4.Range("H4").Value = 8

# Row 5
$q4.Range("A5").Value = 3
Set-TextValue $q4.Range("B5") "001900"
Set-TextValue $q4.Range("C5") "诺安精选价值混合"
Set-TextValue $q4.Range("D5") "0.13"
Set-TextValue $q4.Range("E5") "62.72"
Set-TextValue $q4.Range("F5") "2.07"
Set-TextValue $q4.Range("G5") "0.0027"
$q4.Range("H5").Value = 6

# Row 6
$q4.Range("A6").Value = 4
Set-TextValue $q4.Range("B6") "016259"
Set-TextValue $q4.Range("C6") "鑫元安鑫回报混合C"
Set-TextValue $q4.Range("D6") "0.00"
Set-TextValue $q4.Range("E6") "20.13"
Set-TextValue $q4.Range("F6") "0.97"
$q4.Range("G6").Value = 0
$q4.Range("H6").Value = 8

$total.Select()
